$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.733.49'
$ws.Range("E2").Value = '  -2.48%  '
$ws.Range("D3").Value = '3.805.55'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.41%  '
$ws.Range("D7").Value = '3.805.40'
$ws.Range("E7").Value = '  +1.15%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.74%  '
$ws.Range("E10").Value = '  -4.65%  '
$ws.Range("E11").Value = '  -6.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.466'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000244'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.65%  '
$ws.Range("D15").Value = '4.430.43'
$ws.Range("D16").Value = '3.799.03'
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("D17").Value = '67.703.51'
$ws.Range("E17").Value = '  -2.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.88%  '
$ws.Range("E19").Value = '  -3.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '491.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.742'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.63%  '
$ws.Range("E26").Value = '  +6.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.24'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.38%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  +0.50%  '
$ws.Range("E31").Value = '  -2.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.80'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.50%  '
$ws.Range("E34").Value = '  -3.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.02'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.27%  '
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '461.91'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.329'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.27%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.132'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '49.14'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.45%  '
$ws.Range("E42").Value = '  -2.97%  '
$ws.Range("E43").Value = '  -4.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.76%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.22%  '
$ws.Range("D47").Value = '2.844.11'
$ws.Range("E47").Value = '  -3.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '138.88'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("E49").Value = '  -2.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.69'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.03%  '
